$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-22 13:40:39"
$wsZhCn.Range("G4").Value = "2016-02-22 13:43:03"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-22 13:40:56"
$wsDeDe.Range("G4").Value = "2016-02-22 13:43:33"
